# Update "想去人数" (F column) figures on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> new value mapping for column F
$updates = @{
    2 = 2877
    5 = 6698
    6 = 1611
    7 = 18
    8 = 28
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
